$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$priceUpdates = @{
    "D2"  = "245.90"
    "D3"  = "22.11"
    "D5"  = "0.05866"
    "D6"  = "3.384"
    "D7"  = "6.385"
    "D8"  = "0.8123"
    "D9"  = "0.9615"
    "D10" = "0.1422"
    "D11" = "0.03502"
    "D12" = "0.07336"
    "D13" = "0.03032"
    "D14" = "4.471"
    "D15" = "0.09404"
    "D16" = "0.001596"
    "D17" = "0.04839"
    "D18" = "0.0005893"
    "D19" = "0.006124"
    "D21" = "0.0009821"
    "D22" = "0.00009705"
    "D23" = "3.690"
    "D25" = "0.3254"
    "D40" = "0.03854"
    "D41" = "0.003035"
    "D42" = "0.1074"
    "D43" = "0.002441"
    "D44" = "0.005769"
    "D45" = "0.00005650"
    "D47" = "0.6514"
    "D48" = "0.03622"
}

foreach ($addr in $priceUpdates.Keys) {
    $cell = $ws.Range($addr)
    $cell.NumberFormat = "@"
    $cell.Value = $priceUpdates[$addr]
}

$ws.Range("E18").Value = "17OneONE"
$ws.Range("E48").Value = "47BOLOBOLOWorstin24h"
